$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 112, pushing the existing row 112 (and everything
# below it) down by two rows, so it becomes row 114.
$ws.Rows.Item(112).Resize(2).Insert()

# New row 112: Uva / Red Globe / Primera, Region de O'Higgins, week of 2022-03-08
$ws.Cells.Item(112, 1).Value = 11
$ws.Cells.Item(112, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(112, 3).Value = "Bíobío"
$ws.Cells.Item(112, 4).Value = 44628
$ws.Cells.Item(112, 5).Value = 8
$ws.Cells.Item(112, 6).Value = "Fruta"
$ws.Cells.Item(112, 7).Value = 100109
$ws.Cells.Item(112, 8).Value = "Uva"
$ws.Cells.Item(112, 9).Value = 100109001
$ws.Cells.Item(112, 10).Value = "Uva"
$ws.Cells.Item(112, 11).Value = "Red Globe"
$ws.Cells.Item(112, 12).Value = "Primera"
$ws.Cells.Item(112, 13).Value = 200
$ws.Cells.Item(112, 14).Value = 10000
$ws.Cells.Item(112, 15).Value = 11000
$ws.Cells.Item(112, 16).Value = 10500
$ws.Cells.Item(112, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(112, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(112, 19).Value = 583
$ws.Cells.Item(112, 20).Value = 18

# New row 113: Uva / Red Globe / Segunda, Region de O'Higgins, week of 2022-03-08
$ws.Cells.Item(113, 1).Value = 11
$ws.Cells.Item(113, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(113, 3).Value = "Bíobío"
$ws.Cells.Item(113, 4).Value = 44628
$ws.Cells.Item(113, 5).Value = 8
$ws.Cells.Item(113, 6).Value = "Fruta"
$ws.Cells.Item(113, 7).Value = 100109
$ws.Cells.Item(113, 8).Value = "Uva"
$ws.Cells.Item(113, 9).Value = 100109001
$ws.Cells.Item(113, 10).Value = "Uva"
$ws.Cells.Item(113, 11).Value = "Red Globe"
$ws.Cells.Item(113, 12).Value = "Segunda"
$ws.Cells.Item(113, 13).Value = 100
$ws.Cells.Item(113, 14).Value = 9000
$ws.Cells.Item(113, 15).Value = 9000
$ws.Cells.Item(113, 16).Value = 9000
$ws.Cells.Item(113, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(113, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(113, 19).Value = 500
$ws.Cells.Item(113, 20).Value = 18
